$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "34.308.20"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.797.02"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "226.76"
Set-TextValue "D6" "0.572"
$ws.Range("E6").Value = "  +3.02%  "
$ws.Range("E7").Value = "  +0.02%  "
Set-TextValue "D8" "36.06"
$ws.Range("E8").Value = "  +9.51%  "
Set-TextValue "D9" "0.298"
$ws.Range("E9").Value = "  +1.17%  "
Set-TextValue "D10" "0.0689"
$ws.Range("E10").Value = "  +0.12%  "
Set-TextValue "D11" "0.0960"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").Value = "2.056.81"
$ws.Range("E12").Value = "  +0.17%  "
Set-TextValue "D13" "11.47"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "1.795.77"
$ws.Range("E14").Value = "  +0.14%  "
Set-TextValue "D15" "0.639"
$ws.Range("E15").Value = "  +1.11%  "
Set-TextValue "D16" "4.47"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "34.283.31"
$ws.Range("E17").Value = "  -0.11%  "
Set-TextValue "D18" "68.69"
$ws.Range("E18").Value = "  +0.42%  "
Set-TextValue "D19" "243.54"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("E20").Value = "  -0.59%  "
Set-TextValue "D21" "11.54"
$ws.Range("E21").Value = "  +2.56%  "
$ws.Range("E22").Value = "  +0.04%  "
Set-TextValue "D23" "4.14"
$ws.Range("E23").Value = "  -0.15%  "
Set-TextValue "D24" "2.15"
$ws.Range("E24").Value = "  +3.86%  "
Set-TextValue "D25" "171.97"
$ws.Range("E25").Value = "  +3.38%  "
Set-TextValue "D26" "7.90"
$ws.Range("E26").Value = "  +8.25%  "
Set-TextValue "D27" "16.74"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -0.06%  "
Set-TextValue "D30" "3.98"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D48").Value = "1.957.10"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  +0.02%  "
Set-TextValue "D50" "103.71"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "0.0₆0125"
$ws.Range("E51").Value = "  -2.78%  "

# Row reorderings (coin rank swaps) - update B (name), C (link), D (price), E (volume)
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D31" "0.0526"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D32" "1.24"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D37" "2.44"
$ws.Range("E37").Value = "  -5.99%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D38" "1.06"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "0.955"
$ws.Range("E40").Value = "  +2.12%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D42" "2.81"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D43" "2.41"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D46" "5.99"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D47" "0.0502"
$ws.Range("E47").Value = "  -4.07%  "
